# Weekly update: insert a new price record for "Coliflor" (Vega Monumental
# Concepción) just above the former row 301, shifting every subsequent row
# down by one. The new row mirrors the row above it (row 300) except for a
# new reporting date and a new "Volumen" figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 301 (and everything below it) down by one row.
$ws.Rows(301).Insert()

# Seed the freshly inserted row with the same record as the row above it,
# then overwrite the date and volume with this week's values.
$ws.Range("A300:R300").Copy()
$ws.Range("A301").PasteSpecial()

$ws.Range("D301").Value = 44988
$ws.Range("J301").Value = 1000
$ws.Range("K301").Value = 900
$ws.Range("L301").Value = 1000
$ws.Range("M301").Value = 950
$ws.Range("P301").Value = 950
